$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 107, pushing existing rows 107..204 down to 108..205.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with its data.
$ws.Cells.Item(107, 1).Value = 3
$ws.Cells.Item(107, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(107, 3).Value = "Coquimbo"
$ws.Cells.Item(107, 4).Value = 44484
$ws.Cells.Item(107, 5).Value = 5
$ws.Cells.Item(107, 6).Value = 100114013
$ws.Cells.Item(107, 7).Value = "Zanahoria"
$ws.Cells.Item(107, 8).Value = "Sin especificar"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 160
$ws.Cells.Item(107, 11).Value = 8000
$ws.Cells.Item(107, 12).Value = 8000
$ws.Cells.Item(107, 13).Value = 8000
$ws.Cells.Item(107, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(107, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(107, 16).Value = 400
$ws.Cells.Item(107, 17).Value = 20
$ws.Cells.Item(107, 18).Value = "Hortaliza"
